$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "29.874.98"
Set-TextValue $ws.Range("E2") "  -0.18%  "

Set-TextValue $ws.Range("D3") "1.887.97"
Set-TextValue $ws.Range("E3") "  -0.42%  "

Set-TextValue $ws.Range("D4") "1.000"
Set-TextValue $ws.Range("E4") "  -0.17%  "

Set-TextValue $ws.Range("D5") "0.7494"
Set-TextValue $ws.Range("E5") "  -4.38%  "

Set-TextValue $ws.Range("D6") "242.53"
Set-TextValue $ws.Range("E6") "  -0.65%  "

Set-TextValue $ws.Range("D7") "1.001"
Set-TextValue $ws.Range("E7") "  +0.01%  "

Set-TextValue $ws.Range("D8") "0.3124"
Set-TextValue $ws.Range("E8") "  -0.79%  "

Set-TextValue $ws.Range("D9") "25.47"
Set-TextValue $ws.Range("E9") "  -1.12%  "

Set-TextValue $ws.Range("D10") "0.07136"
Set-TextValue $ws.Range("E10") "  -1.93%  "

Set-TextValue $ws.Range("D11") "0.08491"
Set-TextValue $ws.Range("E11") "  +4.61%  "

Set-TextValue $ws.Range("D12") "0.7606"
Set-TextValue $ws.Range("E12") "  -1.92%  "

Set-TextValue $ws.Range("D13") "1.914.03"
Set-TextValue $ws.Range("E13") "  +3.52%  "

Set-TextValue $ws.Range("D14") "5.362"
Set-TextValue $ws.Range("E14") "  -2.37%  "

Set-TextValue $ws.Range("D15") "93.49"
Set-TextValue $ws.Range("E15") "  -0.92%  "

Set-TextValue $ws.Range("D16") "6.152"
Set-TextValue $ws.Range("E16") "  -1.35%  "

Set-TextValue $ws.Range("D17") "29.965.16"
Set-TextValue $ws.Range("E17") "  +0.05%  "

Set-TextValue $ws.Range("D18") "13.73"
Set-TextValue $ws.Range("E18") "  -1.63%  "

Set-TextValue $ws.Range("D19") "243.61"
Set-TextValue $ws.Range("E19") "  -1.01%  "

Set-TextValue $ws.Range("D20") "0.000007799"
Set-TextValue $ws.Range("E20") "  -0.33%  "

Set-TextValue $ws.Range("D21") "2.162.95"
Set-TextValue $ws.Range("E21") "  +5.24%  "

Set-TextValue $ws.Range("D22") "0.9999"
Set-TextValue $ws.Range("E22") "  +0.04%  "

Set-TextValue $ws.Range("D23") "8.023"
Set-TextValue $ws.Range("E23") "  -1.34%  "

Set-TextValue $ws.Range("E24") "  -0.07%  "

Set-TextValue $ws.Range("D25") "0.1595"
Set-TextValue $ws.Range("E25") "  -0.46%  "

Set-TextValue $ws.Range("D26") "9.385"
Set-TextValue $ws.Range("E26") "  -0.82%  "

Set-TextValue $ws.Range("D27") "162.93"
Set-TextValue $ws.Range("E27") "  -0.49%  "

Set-TextValue $ws.Range("D28") "18.77"
Set-TextValue $ws.Range("E28") "  -0.14%  "

Set-TextValue $ws.Range("D29") "2.030"
Set-TextValue $ws.Range("E29") "  -0.16%  "

Set-TextValue $ws.Range("E30") "  +5.30%  "

Set-TextValue $ws.Range("D31") "1.530"
Set-TextValue $ws.Range("E31") "  -1.19%  "

Set-TextValue $ws.Range("D32") "4.479"
Set-TextValue $ws.Range("E32") "  -0.17%  "

Set-TextValue $ws.Range("D33") "4.107"
Set-TextValue $ws.Range("E33") "  +0.39%  "

Set-TextValue $ws.Range("D34") "0.05402"
Set-TextValue $ws.Range("E34") "  -3.28%  "

Set-TextValue $ws.Range("D35") "1.238"
Set-TextValue $ws.Range("E35") "  -0.75%  "

Set-TextValue $ws.Range("D36") "0.7441"
Set-TextValue $ws.Range("E36") "  -1.27%  "

Set-TextValue $ws.Range("E37") "  +0.53%  "

Set-TextValue $ws.Range("D38") "2.714"
Set-TextValue $ws.Range("E38") "  +1.30%  "

Set-TextValue $ws.Range("D39") "0.01934"
Set-TextValue $ws.Range("E39") "  -0.09%  "

Set-TextValue $ws.Range("D40") "2.769"
Set-TextValue $ws.Range("E40") "  -0.98%  "

Set-TextValue $ws.Range("D41") "0.4460"
Set-TextValue $ws.Range("E41") "  -0.06%  "

Set-TextValue $ws.Range("D42") "1.097.28"
Set-TextValue $ws.Range("E42") "  -3.81%  "

Set-TextValue $ws.Range("D43") "6.083"
Set-TextValue $ws.Range("E43") "  +1.95%  "

Set-TextValue $ws.Range("D44") "72.60"
Set-TextValue $ws.Range("E44") "  -1.66%  "

Set-TextValue $ws.Range("D45") "0.8560"
Set-TextValue $ws.Range("E45") "  +0.10%  "

Set-TextValue $ws.Range("E46") "  +0.02%  "

Set-TextValue $ws.Range("D47") "102.50"
Set-TextValue $ws.Range("E47") "  +0.56%  "

Set-TextValue $ws.Range("D48") "7.668"
Set-TextValue $ws.Range("E48") "  +1.73%  "

Set-TextValue $ws.Range("D49") "1.864"
Set-TextValue $ws.Range("E49") "  -1.79%  "

Set-TextValue $ws.Range("D50") "3.059"
Set-TextValue $ws.Range("E50") "  -2.69%  "

Set-TextValue $ws.Range("D51") "2.058.95"
Set-TextValue $ws.Range("E51") "  +2.35%  "

Write-Output "Updated cryptos list"